$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 40000380
$ws.Range("I33").Value = 58824030
$ws.Range("K33").Value = 58824030
$ws.Range("M33").Value = -58823801

$ws.Range("H40").Value = 1721.5
$ws.Range("I40").Value = 1440.2
$ws.Range("J40").Value = 1877.7778
$ws.Range("K40").Value = 1440.2
$ws.Range("L40").Value = 1877.7778
$ws.Range("M40").Value = -1265.2
$ws.Range("N40").Value = -2227.7778

$ws.Range("H43").Value = 15663
$ws.Range("I43").Value = 5500
$ws.Range("J43").Value = 25826
$ws.Range("K43").Value = 5500
$ws.Range("L43").Value = 25826
$ws.Range("M43").Value = -5431
$ws.Range("N43").Value = -25964

$ws.Range("H139").Value = 50000
$ws.Range("J139").Value = 50000
$ws.Range("L139").Value = 50000
$ws.Range("N139").Value = -60280

$ws.Range("H140").Value = 50200
$ws.Range("J140").Value = 50200
$ws.Range("L140").Value = 50200
$ws.Range("N140").Value = -60560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2700.75
$ws.Range("I45").Value = 2601
$ws.Range("K45").Value = 2601
$ws.Range("M45").Value = -2224

$ws.Range("H58").Value = 39640.25
$ws.Range("J58").Value = 39640.25
$ws.Range("L58").Value = 39640.25
$ws.Range("N58").Value = -40500.25

$ws.Range("H61").Value = 14288560
$ws.Range("I61").Value = 15154245
$ws.Range("K61").Value = 15154245
$ws.Range("M61").Value = -15154033

$ws.Range("H132").Value = 6252209.5
$ws.Range("I132").Value = 7144705
$ws.Range("K132").Value = 21434115
$ws.Range("M132").Value = -21431585

$ws.Range("H136").Value = 14288560
$ws.Range("I136").Value = 15154245
$ws.Range("K136").Value = 45462735
$ws.Range("M136").Value = -45460185

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4457.303
$ws.Range("I105").Value = 2840
$ws.Range("J105").Value = 4892.731
$ws.Range("K105").Value = 2840
$ws.Range("L105").Value = 4892.731
$ws.Range("M105").Value = -1093
$ws.Range("N105").Value = -8386.731

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7940349.5
$ws.Range("I31").Value = 3929.2703
$ws.Range("J31").Value = 66669860
$ws.Range("K31").Value = 3929.2703
$ws.Range("L31").Value = 66669860
$ws.Range("M31").Value = -3634.2703
$ws.Range("N31").Value = -66670450

$ws.Range("H34").Value = 7940349.5
$ws.Range("I34").Value = 3929.2703
$ws.Range("J34").Value = 66669860
$ws.Range("K34").Value = 3929.2703
$ws.Range("L34").Value = 66669860
$ws.Range("M34").Value = -3727.2703
$ws.Range("N34").Value = -66670264

$ws.Range("H58").Value = 3228.625
$ws.Range("I58").Value = 1426.5555
$ws.Range("J58").Value = 5545.5713
$ws.Range("K58").Value = 1426.5555
$ws.Range("L58").Value = 5545.5713
$ws.Range("M58").Value = -1223.5555
$ws.Range("N58").Value = -5951.5713

$ws.Range("H94").Value = 3520.64
$ws.Range("I94").Value = 2162
$ws.Range("J94").Value = 4588.143
$ws.Range("K94").Value = 2162
$ws.Range("L94").Value = 4588.143
$ws.Range("M94").Value = -1711
$ws.Range("N94").Value = -5490.143

$ws.Range("H99").Value = 1212.625
$ws.Range("I99").Value = 1101.2
$ws.Range("J99").Value = 1398.3334
$ws.Range("K99").Value = 1101.2
$ws.Range("L99").Value = 1398.3334
$ws.Range("M99").Value = 396.8
$ws.Range("N99").Value = -4394.3334

$ws.Range("H107").Value = 548.2308
$ws.Range("I107").Value = 444.83334
$ws.Range("J107").Value = 636.8570999999999
$ws.Range("K107").Value = 444.83334
$ws.Range("L107").Value = 636.8570999999999
$ws.Range("M107").Value = 1475.16666
$ws.Range("N107").Value = -4476.8571

$ws.Range("H126").Value = 1212.625
$ws.Range("I126").Value = 1101.2
$ws.Range("J126").Value = 1398.3334
$ws.Range("K126").Value = 3303.6
$ws.Range("L126").Value = 4195.0002
$ws.Range("M126").Value = -833.6000000000004
$ws.Range("N126").Value = -9135.0002

$ws.Range("H136").Value = 3228.625
$ws.Range("I136").Value = 1426.5555
$ws.Range("J136").Value = 5545.5713
$ws.Range("K136").Value = 4279.666499999999
$ws.Range("L136").Value = 16636.7139
$ws.Range("M136").Value = -1729.666499999999
$ws.Range("N136").Value = -21736.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 31.12
$ws.Range("I12").Value = 30.857143
$ws.Range("J12").Value = 31.222221
$ws.Range("K12").Value = 92.57142899999999
$ws.Range("L12").Value = 93.666663
$ws.Range("M12").Value = 80.42857100000001
$ws.Range("N12").Value = -439.666663

$ws.Range("H59").Value = 2333.3333
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 2333.3333
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 6999.999899999999
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -8079.999899999999

$ws.Range("H105").Value = 1666.6666
$ws.Range("J105").Value = 1666.6666
$ws.Range("L105").Value = 4999.9998
$ws.Range("N105").Value = -10241.9998

$ws.Range("H114").Value = 2293.75
$ws.Range("I114").Value = 279
$ws.Range("J114").Value = 2965.3333
$ws.Range("K114").Value = 837
$ws.Range("L114").Value = 8895.999899999999
$ws.Range("M114").Value = 2417
$ws.Range("N114").Value = -15403.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 17953.205
$ws.Range("I70").Value = 26414.186
$ws.Range("J70").Value = 4515.1763
$ws.Range("K70").Value = 26414.186
$ws.Range("L70").Value = 4515.1763
$ws.Range("M70").Value = -26144.186
$ws.Range("N70").Value = -5055.1763

$ws.Range("H73").Value = 17953.205
$ws.Range("I73").Value = 26414.186
$ws.Range("J73").Value = 4515.1763
$ws.Range("K73").Value = 26414.186
$ws.Range("L73").Value = 4515.1763
$ws.Range("M73").Value = -25478.186
$ws.Range("N73").Value = -6387.1763

$ws.Range("H122").Value = 4446261
$ws.Range("I122").Value = 5129762.5
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 15389287.5
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -15386837.5
$ws.Range("N122").Value = -15400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5560
$ws.Range("I7").Value = 5641.6665
$ws.Range("J7").Value = 5437.5
$ws.Range("K7").Value = 5641.6665
$ws.Range("L7").Value = 5437.5
$ws.Range("M7").Value = -5529.6665
$ws.Range("N7").Value = -5661.5

$ws.Range("H22").Value = 1540.75
$ws.Range("I22").Value = 1070
$ws.Range("J22").Value = 1697.6666
$ws.Range("K22").Value = 1070
$ws.Range("L22").Value = 1697.6666
$ws.Range("M22").Value = -775
$ws.Range("N22").Value = -2287.6666

$ws.Range("H27").Value = 1540.75
$ws.Range("I27").Value = 1070
$ws.Range("J27").Value = 1697.6666
$ws.Range("K27").Value = 1070
$ws.Range("L27").Value = 1697.6666
$ws.Range("M27").Value = -963
$ws.Range("N27").Value = -1911.6666

$ws.Range("H122").Value = 11290.4
$ws.Range("I122").Value = 21001.334
$ws.Range("J122").Value = 7128.5713
$ws.Range("K122").Value = 63004.00199999999
$ws.Range("L122").Value = 21385.7139
$ws.Range("M122").Value = -60554.00199999999
$ws.Range("N122").Value = -26285.7139

$ws.Range("H126").Value = 5560
$ws.Range("I126").Value = 5641.6665
$ws.Range("J126").Value = 5437.5
$ws.Range("K126").Value = 16924.9995
$ws.Range("L126").Value = 16312.5
$ws.Range("M126").Value = -14454.9995
$ws.Range("N126").Value = -21252.5

$ws.Range("H132").Value = 9267219
$ws.Range("I132").Value = 5539.2285
$ws.Range("J132").Value = 26328210
$ws.Range("K132").Value = 16617.6855
$ws.Range("L132").Value = 78984630
$ws.Range("M132").Value = -14087.6855
$ws.Range("N132").Value = -78989690

$ws.Range("H139").Value = 55637.5
$ws.Range("J139").Value = 60633.332
$ws.Range("L139").Value = 60633.332
$ws.Range("N139").Value = -70913.33199999999
